# Update cryptocurrency price/volume data per the Sun Jul 2 10:44:01 UTC 2023
# GitHub Actions refresh. Also fixes row 21/22 ordering (Uniswap now ranks
# above Wrapped liquid staked Ether 2.0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @{
    'D2' = '30.513.46'
    'E2' = '  +0.09%  '
    'D3' = '1.916.46'
    'E3' = '  -0.27%  '
    'D4' = '0.9998'
    'E4' = '  -0.07%  '
    'D5' = '245.36'
    'E5' = '  +1.03%  '
    'D6' = '0.9998'
    'E6' = '  -0.07%  '
    'D7' = '0.4796'
    'E7' = '  +2.11%  '
    'D8' = '0.2887'
    'E8' = '  +0.37%  '
    'E9' = '  -0.71%  '
    'D10' = '110.35'
    'E10' = '  +2.83%  '
    'D11' = '19.21'
    'E11' = '  +4.85%  '
    'D12' = '1.912.27'
    'E12' = '  -0.41%  '
    'D13' = '0.07569'
    'E13' = '  -2.23%  '
    'D14' = '5.247'
    'E14' = '  -1.32%  '
    'D15' = '0.6677'
    'E15' = '  +1.40%  '
    'D16' = '302.24'
    'E16' = '  +2.91%  '
    'D17' = '30.515.49'
    'E17' = '  +0.11%  '
    'E18' = '  +0.24%  '
    'D19' = '0.9996'
    'D20' = '0.000007574'
    'E20' = '  -0.53%  '
    'B21' = 'Uniswap'
    'C21' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D21' = '5.506'
    'E21' = '  +4.86%  '
    'B22' = 'WrappedliquidstakedEther2.0'
    'C22' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D22' = '2.164.89'
    'E22' = '  +0.98%  '
    'D23' = '1.0000'
    'E23' = '  -0.19%  '
    'D24' = '6.423'
    'E24' = '  +3.67%  '
    'D25' = '9.474'
    'E25' = '  +1.04%  '
    'D26' = '164.30'
    'E26' = '  -2.70%  '
    'D27' = '20.41'
    'E27' = '  -5.06%  '
    'D28' = '2.105'
    'E28' = '  -0.93%  '
    'D29' = '0.1077'
    'E29' = '  +0.91%  '
    'D30' = '1.399'
    'E30' = '  +2.36%  '
    'D31' = '4.167'
    'E31' = '  -0.46%  '
    'E32' = '  +0.70%  '
    'D33' = '0.04989'
    'E33' = '  -0.91%  '
    'D34' = '0.7366'
    'E34' = '  -0.24%  '
    'D35' = '1.138'
    'E35' = '  -1.36%  '
    'D36' = '0.02050'
    'E36' = '  -2.20%  '
    'D37' = '0.9986'
    'E37' = '  -0.06%  '
    'D38' = '2.725'
    'E38' = '  -0.52%  '
    'E39' = '  -0.35%  '
    'E40' = '  +0.32%  '
    'D41' = '2.020'
    'E41' = '  -2.28%  '
    'D42' = '0.4441'
    'E42' = '  +4.41%  '
    'D43' = '0.8657'
    'E43' = '  -0.64%  '
    'D44' = '5.905'
    'E44' = '  +0.82%  '
    'D45' = '71.15'
    'E45' = '  +5.53%  '
    'E46' = '  -0.04%  '
    'D47' = '50.09'
    'E47' = '  -1.17%  '
    'D48' = '7.282'
    'E48' = '  +1.30%  '
    'D49' = '9.299'
    'E49' = '  -0.01%  '
    'E50' = '  +1.31%  '
    'D51' = '0.2528'
    'E51' = '  +2.55%  '
}

foreach ($addr in $targets.Keys) {
    $rng = $ws.Range($addr)
    # Force the write to land as literal text (matches the original
    # inlineStr cells) instead of letting Excel auto-coerce numeric-looking
    # strings (e.g. "0.9998") into real numbers; then restore the cell's
    # original (default/general) style so no formatting is introduced.
    $rng.NumberFormat = "@"
    $rng.Value = $targets[$addr]
    $rng.Style = "Normal"
}
